$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 71.329076
$ws.Range("H2").Value = 213.987228
$ws.Range("I2").Value = 0.3307464087015077
$ws.Range("J2").Value = 0.3307464087015077
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.16653
$ws.Range("N2").Value = 51.49959
$ws.Range("O2").Value = 0.0560345397128279
$ws.Range("P2").Value = 0.0560345397128279
$ws.Range("Q2").Value = 1224.47272302628
$ws.Range("R2").Value = 11020.25450723652
$ws.Range("S2").Value = 0.01853322277325984
$ws.Range("T2").Value = 0.01853322277325984

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 71.329076
$ws.Range("H3").Value = 213.987228
$ws.Range("I3").Value = 0.3307464087015077
$ws.Range("J3").Value = 0.3307464087015077
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 256.4443053333333
$ws.Range("N3").Value = 769.332916
$ws.Range("O3").Value = 0.8370788162388805
$ws.Range("P3").Value = 0.8370788162388805
$ws.Range("Q3").Value = 18291.93534488854
$ws.Range("R3").Value = 164627.4181039968
$ws.Range("S3").Value = 0.276860812271119
$ws.Range("T3").Value = 0.276860812271119

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 71.329076
$ws.Range("H4").Value = 213.987228
$ws.Range("I4").Value = 0.3307464087015077
$ws.Range("J4").Value = 0.3307464087015077
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.74538866666666
$ws.Range("N4").Value = 98.236166
$ws.Range("O4").Value = 0.1068866440482915
$ws.Range("P4").Value = 0.1068866440482915
$ws.Range("Q4").Value = 2335.698316854205
$ws.Range("R4").Value = 21021.28485168785
$ws.Range("S4").Value = 0.0353523736571288
$ws.Range("T4").Value = 0.03535237365712881

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 129.5615336666667
$ws.Range("H5").Value = 388.684601
$ws.Range("I5").Value = 0.6007649947142101
$ws.Range("J5").Value = 0.6007649947142101
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.16653
$ws.Range("N5").Value = 51.49959
$ws.Range("O5").Value = 0.0560345397128279
$ws.Range("P5").Value = 0.0560345397128279
$ws.Range("Q5").Value = 2224.121954534844
$ws.Range("R5").Value = 20017.09759081359
$ws.Range("S5").Value = 0.03366358995439025
$ws.Range("T5").Value = 0.03366358995439025

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 129.5615336666667
$ws.Range("H6").Value = 388.684601
$ws.Range("I6").Value = 0.6007649947142101
$ws.Range("J6").Value = 0.6007649947142101
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 256.4443053333333
$ws.Range("N6").Value = 769.332916
$ws.Range("O6").Value = 0.8370788162388805
$ws.Range("P6").Value = 0.8370788162388805
$ws.Range("Q6").Value = 33225.31749906962
$ws.Range("R6").Value = 299027.8574916265
$ws.Range("S6").Value = 0.5028876506131282
$ws.Range("T6").Value = 0.5028876506131282

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 129.5615336666667
$ws.Range("H7").Value = 388.684601
$ws.Range("I7").Value = 0.6007649947142101
$ws.Range("J7").Value = 0.6007649947142101
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 32.74538866666666
$ws.Range("N7").Value = 98.236166
$ws.Range("O7").Value = 0.1068866440482915
$ws.Range("P7").Value = 0.1068866440482915
$ws.Range("Q7").Value = 4242.542776164419
$ws.Range("R7").Value = 38182.88498547977
$ws.Range("S7").Value = 0.06421375414669152
$ws.Range("T7").Value = 0.06421375414669152

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.770314
$ws.Range("H8").Value = 44.310942
$ws.Range("I8").Value = 0.06848859658428214
$ws.Range("J8").Value = 0.06848859658428215
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 17.16653
$ws.Range("N8").Value = 51.49959
$ws.Range("O8").Value = 0.0560345397128279
$ws.Range("P8").Value = 0.0560345397128279
$ws.Range("Q8").Value = 253.55503839042
$ws.Range("R8").Value = 2281.99534551378
$ws.Range("S8").Value = 0.003837726985177806
$ws.Range("T8").Value = 0.003837726985177807

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.770314
$ws.Range("H9").Value = 44.310942
$ws.Range("I9").Value = 0.06848859658428214
$ws.Range("J9").Value = 0.06848859658428215
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 256.4443053333333
$ws.Range("N9").Value = 769.332916
$ws.Range("O9").Value = 0.8370788162388805
$ws.Range("P9").Value = 0.8370788162388805
$ws.Range("Q9").Value = 3787.762913285208
$ws.Range("R9").Value = 34089.86621956687
$ws.Range("S9").Value = 0.05733035335463312
$ws.Range("T9").Value = 0.05733035335463314

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.770314
$ws.Range("H10").Value = 44.310942
$ws.Range("I10").Value = 0.06848859658428214
$ws.Range("J10").Value = 0.06848859658428215
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 32.74538866666666
$ws.Range("N10").Value = 98.236166
$ws.Range("O10").Value = 0.1068866440482915
$ws.Range("P10").Value = 0.1068866440482915
$ws.Range("Q10").Value = 483.659672658708
$ws.Range("R10").Value = 4352.937053928373
$ws.Range("S10").Value = 0.007320516244471199
$ws.Range("T10").Value = 0.007320516244471201

Write-Output "done"